$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

$values = @(
    @(11, 50, 300),
    @(13, 46, 300),
    @(11, 22, 300),
    @(11, 6, 300),
    @(11, 0, 300),
    @(11, 1, 300),
    @(13, 10, 300),
    @(12, 0, 300),
    @(11, 12, 300),
    @(11, 39, 300),
    @(12, 34, 300),
    @(10, 34, 300),
    @(11, 29, 300),
    @(11, 47, 300),
    @(11, 31, 300),
    @(10, 48, 300),
    @(12, 28, 300),
    @(11, 35, 300),
    @(10, 54, 300),
    @(11, 25, 300)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i][0]
    $ws.Cells.Item($row, 4).Value = $values[$i][1]
    $ws.Cells.Item($row, 5).Value = $values[$i][2]
}

$ws.Range("D22").Select()
